$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (also updates the workbook-level defined names automatically) ---
$ws.Name = "Level Calibrations"

# --- Update the raw calibration inputs (C2:E3) ---
$ws.Range("C2").Value = 8.5
$ws.Range("D2").Value = 1.1195999999999999
$ws.Range("E2").Value = 8.4

$ws.Range("C3").Value = 1.956
$ws.Range("D3").Value = 7.6657999999999999
$ws.Range("E3").Value = 2

# --- D12 / D13 formulas change to reference the Scaled column (C) instead of Keyence (D) ---
$ws.Range("D12").Formula = "=(E3-E2)/(C3-C2)"
$ws.Range("C13").Value = "Lse"
$ws.Range("D13").Formula = "=C2"

# --- New numeric format used by the new table (rows 20-27) ---
$fmt = "0.0000"

# --- New header row (19) for the Key Measure calibration-check table ---
$ws.Range("B19").Value = "Key Measure"
$ws.Range("C19").Value = "Corrected"
$ws.Range("C19").WrapText = $true
$ws.Range("D19").Value = "Scaled"
$ws.Range("D19").WrapText = $true

# --- Row 20 (Empty anchor row) ---
$ws.Range("A20").Value = "Empty"
$ws.Range("B20").Value = 1
$ws.Range("C20").Formula = "=(D20-`$D`$13)*`$D`$12+`$D`$14"
$ws.Range("D20").Formula = "=B20*`$D`$5+`$D`$6"
$ws.Range("C20").NumberFormat = $fmt
$ws.Range("C20").WrapText = $true
$ws.Range("D20").NumberFormat = $fmt
$ws.Range("D20").WrapText = $true

# --- Row 21 ---
$ws.Range("B21").Value = 1.1195999999999999
$ws.Range("C21:C27").Formula = "=(D21-`$D`$13)*`$D`$12+`$D`$14"
$ws.Range("D21:D27").Formula = "=B21*`$D`$5+`$D`$6"
$ws.Range("C21:D27").NumberFormat = $fmt
$ws.Range("C21:D27").WrapText = $true

# --- Rows 22-25 (plain B inputs) ---
$ws.Range("B22").Value = 3
$ws.Range("B23").Value = 4
$ws.Range("B24").Value = 5
$ws.Range("B25").Value = 6

# --- Row 26 ---
$ws.Range("B26").Value = 7.6657999999999999

# --- Row 27 (Full anchor row) ---
$ws.Range("A27").Value = "Full"
$ws.Range("B27").Value = 8

# --- Selection state to match the saved workbook ---
$ws.Range("C2:D3").Select()
